# Update cryptocurrency symbol list (prices + reordered rows) per Dec 16 2022 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    # Preserve these numeric-looking values as literal text (matches the inlineStr
    # cells in the workbook) instead of letting Excel coerce them to numbers.
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range('D2') '261.22'
# Row 3
Set-TextValue $ws.Range('D3') '22.88'
# Row 4
Set-TextValue $ws.Range('D4') '6.185'
# Row 5
Set-TextValue $ws.Range('D5') '0.06235'
# Row 6
Set-TextValue $ws.Range('D6') '6.733'
# Row 7
Set-TextValue $ws.Range('D7') '3.444'
# Row 8
Set-TextValue $ws.Range('D8') '1.345'
# Row 9
Set-TextValue $ws.Range('D9') '0.7967'
# Row 10
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range('D10') '0.1572'
$ws.Range('E10').Value = '9WazirXWRX'
# Row 11
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range('D11') '0.08097'
$ws.Range('E11').Value = '10MandalaExchangeTokenMDX'
# Row 12
$ws.Range('B12').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C12').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range('D12') '0.03437'
$ws.Range('E12').Value = '11LiechtensteinCryptoassetsExchangeLCX'
# Row 13
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range('D13') '0.03086'
$ws.Range('E13').Value = '12BitrueCoinBTR'
# Row 14
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range('D14') '0.09324'
$ws.Range('E14').Value = '13BitMartTokenBMX'
# Row 15
$ws.Range('B15').Value = 'MCDex'
$ws.Range('C15').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws.Range('D15') '3.847'
$ws.Range('E15').Value = '14MCDexMCB'
# Row 16
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range('D16') '0.001689'
$ws.Range('E16').Value = '15BitForexTokenBF'
# Row 17
$ws.Range('B17').Value = 'CoinExToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws.Range('D17') '0.04796'
$ws.Range('E17').Value = '16CoinExTokenCET'
# Row 18
$ws.Range('B18').Value = 'One'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws.Range('D18') '0.0006148'
$ws.Range('E18').Value = '17OneONEWorstin24h'
# Row 19
Set-TextValue $ws.Range('D19') '0.006211'
# Row 20
Set-TextValue $ws.Range('D20') '0.006159'
# Row 21
Set-TextValue $ws.Range('D21') '0.001092'
# Row 22
Set-TextValue $ws.Range('D22') '0.0001499'
# Row 23
Set-TextValue $ws.Range('D23') '3.697'
# Row 24
Set-TextValue $ws.Range('D24') '2.207'
# Row 25
Set-TextValue $ws.Range('D25') '0.3339'
# Row 40
Set-TextValue $ws.Range('D40') '0.04612'
# Row 41
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue $ws.Range('D41') '0.1118'
$ws.Range('E41').Value = '40BKEXTokenBKK'
# Row 42
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue $ws.Range('D42') '0.003129'
$ws.Range('E42').Value = '41CEJICEJI'
# Row 43
$ws.Range('B43').Value = 'KickToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue $ws.Range('D43') '0.007077'
$ws.Range('E43').Value = '42KickTokenKICK'
# Row 44
Set-TextValue $ws.Range('D44') '0.01011'
# Row 46
Set-TextValue $ws.Range('D46') '0.00005879'
# Row 49
Set-TextValue $ws.Range('D49') '0.09020'
